# Clear the "Name" property value on the "Test" sheet (B2), then make
# "Test" the active/selected sheet with B3 as the selected cell -- this
# mirrors dropping the stale "Karr et al., 2016 ..." fixture string and
# switching the workbook's active tab from "Property" (index 1) back to
# "Test" (index 0) with cell B3 selected.

$wb = $excel.ActiveWorkbook

$testSheet = $wb.Worksheets.Item("Test")

# Remove the old "Name" value (row 2, column B) entirely.
$testSheet.Range("B2").Value = ""

# "Test" becomes the active sheet / tab, with B3 selected.
$testSheet.Activate()
$testSheet.Range("B3").Select()
